# Refresh the cryptocurrency price/volume snapshot (columns D "Price" and
# E "Volume(1h)") for the rows whose figures changed in this update run.
#
# Some "Price" values look like plain numbers (e.g. "22.60", "1.001") but
# must stay exactly as text (matching the original inlineStr cells, no
# rounding/trailing-zero loss). Prefixing with a single quote makes Excel
# treat the entry as text (quote-prefixed), just like typing '22.60 into
# the cell would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.072.22'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.836.22'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''243.21'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").Value = '''0.6286'
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.07615'
$ws.Range("E8").Value = '  +3.68%  '
$ws.Range("D9").Value = '''0.2934'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").Value = '''22.66'
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("D11").Value = '''0.07744'
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("D12").Value = '1.848.51'
$ws.Range("E12").Value = '  +0.94%  '
$ws.Range("D13").Value = '''4.969'
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").Value = '''0.6663'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("E15").Value = '  +16.57%  '
$ws.Range("D16").Value = '''82.91'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").Value = '''6.068'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '29.092.29'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = '''227.73'
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").Value = '''12.40'
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").Value = '''7.230'
$ws.Range("E22").Value = '  +1.57%  '
$ws.Range("D24").Value = '''159.40'
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("D25").Value = '''8.526'
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("D26").Value = '''0.1385'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").Value = '''17.95'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("D28").Value = '''1.494'
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").Value = '''4.108'
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("D30").Value = '''4.026'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '''1.194'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").Value = '''0.05263'
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("D33").Value = '''1.844'
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("D34").Value = '''0.7357'
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").Value = '''1.140'
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").Value = '''2.698'
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("D37").Value = '1.244.14'
$ws.Range("E37").Value = '  -3.43%  '
$ws.Range("D38").Value = '''2.761'
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").Value = '''0.01789'
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").Value = '''6.370'
$ws.Range("E40").Value = '  +0.42%  '
$ws.Range("D41").Value = '''0.8974'
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D43").Value = '''102.12'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").Value = '1.986.99'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = '''0.00000000123'
$ws.Range("E45").Value = '  +3.13%  '
$ws.Range("D46").Value = '''64.48'
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("D47").Value = '''0.5114'
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").Value = '''0.4048'
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("D49").Value = '''8.903'
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("D50").Value = '''0.05759'
$ws.Range("E50").Value = '  -1.16%  '
$ws.Range("D51").Value = '''6.697'
$ws.Range("E51").Value = '  +0.18%  '
